$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.120.73"
$ws.Range("E2").Value = "  +3.53%  "
$ws.Range("D3").Value = "1.723.57"
$ws.Range("E3").Value = "  +2.55%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "'218.95"
$ws.Range("E5").Value = "  +1.71%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("D8").Value = "'24.30"
$ws.Range("E8").Value = "  +13.60%  "
$ws.Range("E9").Value = "  +3.19%  "
$ws.Range("E10").Value = "  +1.53%  "
$ws.Range("D11").Value = "'0.0898"
$ws.Range("E11").Value = "  +1.41%  "
$ws.Range("D12").Value = "1.967.68"
$ws.Range("E12").Value = "  +2.63%  "
$ws.Range("D13").Value = "1.730.81"
$ws.Range("E13").Value = "  +2.94%  "
$ws.Range("E14").Value = "  +3.05%  "
$ws.Range("E15").Value = "  +4.74%  "
$ws.Range("E16").Value = "  +1.97%  "
$ws.Range("D17").Value = "28.082.71"
$ws.Range("E17").Value = "  +3.42%  "
$ws.Range("D18").Value = "'243.01"
$ws.Range("E18").Value = "  +1.43%  "
$ws.Range("D19").Value = "'8.03"
$ws.Range("E19").Value = "  -0.62%  "
$ws.Range("E20").Value = "  +1.25%  "
$ws.Range("E21").Value = "  -0.08%  "
$ws.Range("D22").Value = "'4.62"
$ws.Range("E22").Value = "  +2.25%  "
$ws.Range("E23").Value = "  +2.06%  "
$ws.Range("D24").Value = "'2.12"
$ws.Range("E24").Value = "  -0.16%  "
$ws.Range("D25").Value = "'149.16"
$ws.Range("E25").Value = "  +1.53%  "
$ws.Range("E26").Value = "  +3.53%  "
$ws.Range("E27").Value = "  +2.00%  "
$ws.Range("E28").Value = "  +0.64%  "
$ws.Range("E29").Value = "  -0.26%  "
$ws.Range("E30").Value = "  +2.07%  "
$ws.Range("E31").Value = "  +1.62%  "
$ws.Range("D32").Value = "'3.44"
$ws.Range("E32").Value = "  +2.21%  "
$ws.Range("D33").Value = "1.497.13"
$ws.Range("E33").Value = "  -4.06%  "
$ws.Range("E34").Value = "  +1.74%  "
$ws.Range("E35").Value = "  -1.54%  "
$ws.Range("D36").Value = "'0.959"
$ws.Range("E36").Value = "  +2.59%  "
$ws.Range("D37").Value = "'0.608"
$ws.Range("E37").Value = "  +0.81%  "
$ws.Range("E38").Value = "  +0.78%  "
$ws.Range("E39").Value = "  +0.28%  "
$ws.Range("E40").Value = "  +1.12%  "
$ws.Range("E41").Value = "  +2.12%  "
$ws.Range("D42").Value = "'5.80"
$ws.Range("E42").Value = "  +2.29%  "
$ws.Range("E44").Value = "  +1.66%  "
$ws.Range("D45").Value = "1.872.35"
$ws.Range("E45").Value = "  +2.50%  "
$ws.Range("D46").Value = "'0.805"
$ws.Range("E46").Value = "  +3.07%  "
$ws.Range("E47").Value = "  +10.69%  "
$ws.Range("D48").Value = "'90.76"
$ws.Range("E48").Value = "  +0.03%  "
$ws.Range("E49").Value = "  +5.89%  "
$ws.Range("D50").Value = "'8.25"
$ws.Range("E50").Value = "  +1.09%  "
$ws.Range("E51").Value = "  +0.51%  "

# Normalize formatting: the value assignments above may cause Excel to
# auto-apply a text/quote-prefix number format on cells whose new text
# looks numeric (e.g. "24.30", "0.0898"). The source cells carry no
# explicit style, so clear any such auto-formatting to keep the cells
# styled exactly as before (default style, plain text).
$ws.Range("D2:E51").ClearFormats()
